$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Paragraphs.Item(1).Range.Text = "2024-10-29 Tuesday"

# Update the division problems in the table. Each "data" row of the table
# is followed by three blank rows, so the populated rows are 1, 5, 9, 13, 17.
$tbl = $d.Tables.Item(1)

$updates = @{
    1  = @("93÷5=", "85÷8=", "31÷4=", "52÷5=", "57÷8=")
    5  = @("31÷7=", "25÷6=", "36÷4=", "87÷6=", "28÷2=")
    9  = @("37÷7=", "70÷7=", "85÷3=", "60÷5=", "65÷7=")
    13 = @("86÷2=", "79÷6=", "12÷3=", "58÷8=", "37÷6=")
    17 = @("81÷2=", "34÷8=", "61÷3=", "37÷9=", "77÷7=")
}

foreach ($rowIndex in $updates.Keys) {
    $values = $updates[$rowIndex]
    for ($col = 1; $col -le $values.Length; $col++) {
        $cell = $tbl.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
